# EASYCIVILCOMMANDS.xlsx update
#
# Content changes:
#  - "Coordination" sheet: remove the "DBN / [DBN]..Distance Between All Pipes" row
#  - "Profiles" sheet: remove the "ADP / [ADP]..Auto Add Parts to Profile" row and the
#    "LPP / [LPP]..Layout Profile" row (plus its following separator row)
#  - "Update" sheet: rename the update-button caption from
#    "[ Update EASYCIVIL ]" to "[ Check Update EASYCIVIL ]"
#
# View/selection changes:
#  - "Profiles" sheet: scroll back to top, select B5
#  - "Update" sheet: select B8, no longer the active tab
#  - "Coordination" sheet: select B8, becomes the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Coordination sheet: delete the DBN row
# ---------------------------------------------------------------------------
$wsCoord = $wb.Worksheets.Item("Coordination")
$wsCoord.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# Profiles sheet: delete the ADP row, and the LPP + following separator rows
# (delete from the bottom up so earlier row numbers stay valid)
# ---------------------------------------------------------------------------
$wsProfiles = $wb.Worksheets.Item("Profiles")
$wsProfiles.Rows.Item(11).Delete()
$wsProfiles.Rows.Item(10).Delete()
$wsProfiles.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# Update sheet: change the button caption text
# ---------------------------------------------------------------------------
$wsUpdate = $wb.Worksheets.Item("Update")
$wsUpdate.Range("B1").Value = "[ Check Update EASYCIVIL ]"

# ---------------------------------------------------------------------------
# View state: update selections / active sheet to match the saved workbook
# ---------------------------------------------------------------------------
$wsProfiles.Activate()
$wsProfiles.Range("B5").Select()

$wsUpdate.Activate()
$wsUpdate.Range("B8").Select()

$wsCoord.Activate()
$wsCoord.Range("B8").Select()
